# Understand how to store and return multiple values using a POJO class
# in a return data class.
#
# This adds a new "Dynamic" worksheet (after the existing "data1" sheet)
# that lists the "Name" column header plus a newly added name, "Sumit".

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the last sheet in the workbook (it becomes the
# active / selected sheet, which also clears the previous "tabSelected"
# flag on the former last sheet automatically).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Dynamic"

# Header cell, styled the same way the other "Name" header cells in the
# workbook are (solid fill), followed by the new row of data.
$ws.Range("A1").Value = "Name"
$ws.Range("A1").Interior.Color = 9868950
$ws.Range("A2").Value = "Sumit"

# Size the column to fit its contents, like the other data sheets.
$ws.Columns.Item(1).AutoFit() | Out-Null
